# Adds the 2022-Q3 quarterly data:
#   - a new "2022-Q3" worksheet (inserted right before the existing "2022-Q2"
#     worksheet, cloned from it so it keeps identical formatting/structure)
#   - a new leading data row on the "总计" (totals) summary worksheet
#
# Helper: write a value into a cell as literal TEXT (not auto-converted to a
# number by the COM layer) and strip any stray number-format style picked up
# along the way, so the cell ends up styled exactly like its neighbours.
function Set-TextCell($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Clone the current "2022-Q2" sheet (position 2) into a new sheet placed
#    right before it, then rename the clone to "2022-Q3" and fill in the new
#    quarter's numbers. Every other existing sheet is left untouched - it
#    simply shifts one tab to the right.
# ---------------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item(2)
$q2Sheet.Copy($q2Sheet, $null)
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

$q3Rows = @(
    @("000727", "融通健康产业灵活配置混合A", "22.64", "93.67", "5.73", "1.2973", 6),
    @("009274", "融通健康产业灵活配置混合C", "17.64", "93.67", "5.73", "1.0108", 6),
    @("009805", "国泰医药健康股票A",         "8.71",  "94.26", "6.02", "0.5243", 8),
    @("011404", "融通鑫新成长混合C",         "1.75",  "94.07", "9.83", "0.1720", 1),
    @("011403", "融通鑫新成长混合A",         "0.39",  "94.07", "9.83", "0.0383", 1),
    @("011326", "国泰医药健康股票C",         "0.52",  "94.26", "6.02", "0.0313", 8)
)

$r = 2
foreach ($row in $q3Rows) {
    Set-TextCell $q3Sheet.Range("B$r") $row[0]
    Set-TextCell $q3Sheet.Range("C$r") $row[1]
    Set-TextCell $q3Sheet.Range("D$r") $row[2]
    Set-TextCell $q3Sheet.Range("E$r") $row[3]
    Set-TextCell $q3Sheet.Range("F$r") $row[4]
    Set-TextCell $q3Sheet.Range("G$r") $row[5]
    $q3Sheet.Range("H$r").Value = $row[6]
    $r++
}

# Restore the originally-selected tab (cloning the sheet made the new clone
# the active one; the real active tab is the last sheet, "2021-Q2").
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()

# ---------------------------------------------------------------------------
# 2) Insert the new leading row on the "总计" summary sheet and carry the
#    formatting down from the row that is about to become row 3, so the new
#    row 2 ends up styled identically (column A keeps its centered style,
#    columns B:D stay plain/unstyled).
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$totalSheet.Rows(2).Insert()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B3:D3").Copy()
$totalSheet.Range("B2:D2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 6
$totalSheet.Range("D2").Value = 3.07
